# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from serial date 45224 (2023-10-25) to serial date 45233 (2023-11-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45233
}
